$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.122.15"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "2.485.92"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'323.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'106.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.89%  "
$ws.Range("D7").Value = "'0.526"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "'38.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.39%  "
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "'18.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "2.876.78"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").Value = "2.475.24"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "47.066.01"
$ws.Range("E18").Value = "  +4.08%  "
$ws.Range("D19").Value = "'12.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("D20").Value = "'6.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.15%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value = "'70.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'2.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.61%  "
$ws.Range("D24").Value = "'250.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("E25").Value = "  +4.13%  "
$ws.Range("D26").Value = "'26.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("E29").Value = "  +3.92%  "
$ws.Range("E30").Value = "  +6.28%  "
$ws.Range("E31").Value = "  +9.46%  "
$ws.Range("D32").Value = "'49.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'19.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("E34").Value = "  +4.85%  "
$ws.Range("D35").Value = "'0.0794"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.49%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  +6.25%  "
$ws.Range("D38").Value = "'4.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.15%  "
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("D41").Value = "'122.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "'21.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "1.969.19"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").Value = "'2.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").Value = "'1.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("E50").Value = "  +9.45%  "
$ws.Range("D51").Value = "'79.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.86%  "
